$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max") - shifts "prediction" and "rejection-f" left
$ws.Range("C:C").Delete()

# Delete row 4 (even_MAG-GUT49487.fa) - shifts subsequent rows up
$ws.Range("4:4").Delete()

# Update the B column (now "1-s__Proteus mirabilis") values for data rows
$ws.Cells.Item(2, 2).Value = 0.01796816099174947
$ws.Cells.Item(3, 2).Value = 0.4408786776011882
$ws.Cells.Item(4, 2).Value = -0.01643298445229746
$ws.Cells.Item(5, 2).Value = 0.03508223360167051

# Update column C (now "prediction") - was numeric "1", needs to become text "s__Proteus mirabilis"
$ws.Cells.Item(2, 3).Value = "s__Proteus mirabilis"
$ws.Cells.Item(3, 3).Value = "s__Proteus mirabilis"
$ws.Cells.Item(4, 3).Value = "s__Proteus mirabilis"
$ws.Cells.Item(5, 3).Value = "s__Proteus mirabilis"

# Row 4 (rejection-f, column D) becomes "(reject)" suffixed
$ws.Cells.Item(4, 4).Value = "s__Proteus mirabilis(reject)"
